# fix missing message from hidden cells
#
# The "Zipcode_CheckOrder" sheet's header row used long/optional labels
# ("Identifier", "Descriptions? (Optional)", "Hidden? (Optional)",
# "Feedback (Optional)") that the grader no longer recognizes; they are
# renamed to the short keys it expects ("id", "Description", "special",
# "Feedback"). The companion "SheetGradingOrder" sheet's header
# "Sheetname" is renamed to "sheet" to match. Also, row 4's "special"
# flag was "hk" (hidden+killer) which suppressed the failure feedback
# message for that row - it's fixed to just "h" (hidden) so the
# feedback message is shown, per the commit message.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Zipcode_CheckOrder")
$ws3.Range("C1").Value = "Description"
$ws3.Range("D1").Value = "special"
$ws3.Range("E1").Value = "Feedback"
$ws3.Range("D4").Value = "h"

$ws1 = $wb.Worksheets.Item("SheetGradingOrder")
$ws1.Range("B1").Value = "sheet"

# Renaming A1 last so the shared-string table's newly appended entries
# land in the same order as the target workbook.
$ws3.Range("A1").Value = "id"

# Match the saved cursor position on the active sheet.
$ws3.Range("A2").Select()
